$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four cells in row 2 to the new "16/11/2023" set of values,
# mirroring the previous "15/11/2023" entries that were already present
# (same columns: AP_Value, DT_Value, VR_Value, UF_Prefix).
$ws.Range("N2").Value  = "Appointment Date : 16/11/2023, Time : [ 09:05 AM to 09:09 AM ]"
$ws.Range("AB2").Value = "16/11/2023"
$ws.Range("AR2").Value = "voice_record_16112023"
$ws.Range("AU2").Value = "formshow_16112023"
